$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("D4").Value = 0.13
$ws.Range("F4").Value = -1

$ws.Range("D6").Value = 0.16
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1

$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0.23
$ws.Range("E12").Value = -1

$ws.Range("D14").Value = 0.25

$ws.Range("D15").Value = 0.23
$ws.Range("H15").Value = 0.1
$ws.Range("J15").Value = -100

# Update the shared formula coefficient for H column from *2 to *5, across N3:N15
$ws.Range("N3:N14").Formula = "=C3-D3*20-E3*0.8-F3*0.6-H3*5+J3/300"
$ws.Range("N15").Formula = "=C15-D15*20-E15*0.8-F15*0.6-H15*5+J15/300"

# Update selection to E12
$ws.Range("E12").Select() | Out-Null

# Update column B width
$ws.Columns.Item(2).ColumnWidth = 39.140625
